# Refresh the HidroElektrana log sheet with the latest simulation run.
# The underlying app reran the hydro-plant simulation (this time with the
# new "don't let Use% pass 100" guard) and logged its ticks again, which
# produced fresh timestamps/percentages and a handful of extra rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$times = @(
    "09/01/2023 11:16:03 pm","09/01/2023 11:16:04 pm","09/01/2023 11:16:04 pm","09/01/2023 11:16:04 pm",
    "09/01/2023 11:16:04 pm","09/01/2023 11:16:04 pm","09/01/2023 11:16:05 pm","09/01/2023 11:16:05 pm",
    "09/01/2023 11:16:06 pm","09/01/2023 11:16:06 pm","09/01/2023 11:16:06 pm","09/01/2023 11:16:07 pm",
    "09/01/2023 11:16:07 pm","09/01/2023 11:16:08 pm","09/01/2023 11:16:08 pm","09/01/2023 11:16:09 pm",
    "09/01/2023 11:16:10 pm","09/01/2023 11:16:10 pm","09/01/2023 11:16:10 pm","09/01/2023 11:16:11 pm",
    "09/01/2023 11:16:12 pm","09/01/2023 11:16:12 pm","09/01/2023 11:16:13 pm","09/01/2023 11:16:13 pm",
    "09/01/2023 11:16:14 pm","09/01/2023 11:16:14 pm","09/01/2023 11:16:15 pm","09/01/2023 11:16:15 pm",
    "09/01/2023 11:16:16 pm","09/01/2023 11:16:16 pm","09/01/2023 11:16:17 pm"
)

$usePct = @(0,0,0,0,0,0,1,4,7,10,13,16,19,22,25,28,31,34,37,40,43,46,49,52,55,58,61,64,67,70,73)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $no = ($i + 1).ToString()

    # Column A ("No.") and column B ("Time") are logged as plain text in
    # this sheet, even though "No." looks numeric - force text storage
    # (matching the sheet's existing General-styled text cells) and then
    # drop back to the Normal style so no stray numeric formatting sticks.
    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $no
    $ws.Range("A$row").Style = "Normal"

    $ws.Range("B$row").Value = $times[$i]

    $ws.Range("C$row").Value = $usePct[$i]
}

$ws.Columns.Item(2).AutoFit()
